# Update FFXIV market-board price/profit figures (currentAveragePrice*, LevePrice*, LeveProfit*)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve-profit tables, per scheduled data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 11: Gotta Bounce / Rubber
$ws.Range("H11").Value = 33
$ws.Range("I11").Value = 33
$ws.Range("K11").Value = 33
$ws.Range("M11").Value = 107

# Row 88: The Grave of Hemlock Groves / Growth Formula Zeta
$ws.Range("H88").Value = 1120
$ws.Range("I88").Value = 1600
$ws.Range("K88").Value = 1600
$ws.Range("M88").Value = -1194

# Row 91: Dappling the Highlands (L) / Growth Formula Zeta
$ws.Range("H91").Value = 1120
$ws.Range("I91").Value = 1600
$ws.Range("K91").Value = 1600
$ws.Range("M91").Value = -196

# Row 99: Rumor Has It / Commanding Craftsman's Tea
$ws.Range("H99").Value = 416.33334
$ws.Range("J99").Value = 200
$ws.Range("L99").Value = 600
$ws.Range("N99").Value = -3596

# Row 138: All-night Crafting / Cunning Craftsman's Tisane
$ws.Range("H138").Value = 5120.9644
$ws.Range("I138").Value = 4454.1113
$ws.Range("J138").Value = 5436.8423
$ws.Range("K138").Value = 13362.3339
$ws.Range("L138").Value = 16310.5269
$ws.Range("M138").Value = -8222.333899999998
$ws.Range("N138").Value = -26590.5269

$ws = $wb.Worksheets.Item("ARM")
# Row 5: The Alloyed Truth / Bronze Rivets
$ws.Range("H5").Value = 423.5
$ws.Range("I5").Value = 231.33333
$ws.Range("K5").Value = 231.33333
$ws.Range("M5").Value = -119.33333

# Row 46: Get Me the Usual / Heavy Steel Flanchard
$ws.Range("H46").Value = 9708
$ws.Range("J46").Value = 9694
$ws.Range("L46").Value = 9694
$ws.Range("N46").Value = -10332

# Row 122: Haste for High Durium / High Durium Nugget
$ws.Range("H122").Value = 52006
$ws.Range("I122").Value = 52006
$ws.Range("K122").Value = 156018
$ws.Range("M122").Value = -153568

# Row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Range("H132").Value = 2614.5
$ws.Range("I132").Value = 2614.5
$ws.Range("K132").Value = 7843.5
$ws.Range("M132").Value = -5313.5

$ws = $wb.Worksheets.Item("BSM")
# Row 4: Mending Fences / Bronze Rivets
$ws.Range("H4").Value = 423.5
$ws.Range("I4").Value = 231.33333
$ws.Range("K4").Value = 231.33333
$ws.Range("M4").Value = -116.33333

# Row 86: Through Thick and Thin / Adamantite Nugget
$ws.Range("H86").Value = 1263
$ws.Range("I86").Value = 1056.8334
$ws.Range("K86").Value = 1056.8334
$ws.Range("M86").Value = 66.16660000000002

# Row 89: Piercing Eyes Deserve Piercing Shafts (L) / Adamantite Nugget
$ws.Range("H89").Value = 1263
$ws.Range("I89").Value = 1056.8334
$ws.Range("K89").Value = 5284.166999999999
$ws.Range("M89").Value = 331.8330000000005

# Row 106: Fire for Hire / Molybdenum Rimfire
$ws.Range("H106").Value = 7622.5
$ws.Range("J106").Value = 7622.5
$ws.Range("L106").Value = 7622.5
$ws.Range("N106").Value = -10146.5

# Row 124: History of the Hrothgar / High Durium Bayonet
$ws.Range("H124").Value = 85945
$ws.Range("J124").Value = 85945
$ws.Range("L124").Value = 85945
$ws.Range("N124").Value = -95765

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found / Walnut Lumber
$ws.Range("H31").Value = 1719.1786
$ws.Range("I31").Value = 1274.6
$ws.Range("K31").Value = 1274.6
$ws.Range("M31").Value = -979.5999999999999

# Row 34: Armoires of the Rich and Famous / Walnut Lumber
$ws.Range("H34").Value = 1719.1786
$ws.Range("I34").Value = 1274.6
$ws.Range("K34").Value = 1274.6
$ws.Range("M34").Value = -1072.6

# Row 58: You Do the Heavy Lifting / Mahogany Lumber
$ws.Range("H58").Value = 3391
$ws.Range("I58").Value = 2740.75
$ws.Range("K58").Value = 2740.75
$ws.Range("M58").Value = -2537.75

# Row 62: Splinter in the Sewers / Cedar Lumber
$ws.Range("H62").Value = 1999.5
$ws.Range("I62").Value = 1999.5
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 1999.5
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -1375.5
$ws.Range("N62").ClearContents()

# Row 65: The Lumber of Their Discontent (L) / Cedar Lumber
$ws.Range("H65").Value = 1999.5
$ws.Range("I65").Value = 1999.5
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 9997.5
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -6877.5
$ws.Range("N65").ClearContents()

# Row 88: Hold on Adamantite / Adamantite Spear
$ws.Range("H88").Value = 23225.25
$ws.Range("J88").Value = 23225.25
$ws.Range("L88").Value = 23225.25
$ws.Range("N88").Value = -24037.25

# Row 91: Spears for Stone Vigilantes (L) / Adamantite Spear
$ws.Range("H91").Value = 23225.25
$ws.Range("J91").Value = 23225.25
$ws.Range("L91").Value = 23225.25
$ws.Range("N91").Value = -26033.25

# Row 99: O Pine / Pine Lumber
$ws.Range("H99").Value = 1835.4
$ws.Range("I99").Value = 1542.6666
$ws.Range("K99").Value = 1542.6666
$ws.Range("M99").Value = -44.66660000000002

# Row 108: Just Starting Out / White Oak Fishing Rod
$ws.Range("H108").Value = 62824.5
$ws.Range("J108").Value = 62824.5
$ws.Range("L108").Value = 62824.5
$ws.Range("N108").Value = -70504.5

# Row 126: A Better Conductor / Red Pine Lumber
$ws.Range("H126").Value = 1835.4
$ws.Range("I126").Value = 1542.6666
$ws.Range("K126").Value = 4627.9998
$ws.Range("M126").Value = -2157.9998

# Row 136: Turali Quality / Dark Mahogany Lumber
$ws.Range("H136").Value = 3391
$ws.Range("I136").Value = 2740.75
$ws.Range("K136").Value = 8222.25
$ws.Range("M136").Value = -5672.25

$ws = $wb.Worksheets.Item("CUL")
# Row 12: Butter Me Up / Kukuru Butter
$ws.Range("H12").Value = 59.25
$ws.Range("J12").Value = 59.25
$ws.Range("L12").Value = 177.75
$ws.Range("N12").Value = -523.75

# Row 76: Old Victories, New Tastes / Dhalmel Fricassee
$ws.Range("H76").Value = 13112.556
$ws.Range("I76").Value = 6753.25
$ws.Range("K76").Value = 20259.75
$ws.Range("M76").Value = -19876.75

# Row 79: The Eats of Authenticity (L) / Dhalmel Fricassee
$ws.Range("H79").Value = 13112.556
$ws.Range("I79").Value = 6753.25
$ws.Range("K79").Value = 20259.75
$ws.Range("M79").Value = -18933.75

# Row 80: Saucy for a Suitor / Hollandaise Sauce
$ws.Range("H80").Value = 2988.5
$ws.Range("J80").Value = 2988.5
$ws.Range("L80").Value = 8965.5
$ws.Range("N80").Value = -10837.5

# Row 83: Saved by the Sauce (L) / Hollandaise Sauce
$ws.Range("H83").Value = 2988.5
$ws.Range("J83").Value = 2988.5
$ws.Range("L83").Value = 26896.5
$ws.Range("N83").Value = -36256.5

$ws = $wb.Worksheets.Item("GSM")
# Row 122: Awarding Academic Excellence / Ametrine
$ws.Range("H122").Value = 3492.1667
$ws.Range("I122").Value = 3513.2
$ws.Range("J122").Value = 3387
$ws.Range("K122").Value = 10539.6
$ws.Range("L122").Value = 10161
$ws.Range("M122").Value = -8089.599999999999
$ws.Range("N122").Value = -15061

$ws = $wb.Worksheets.Item("LTW")
# Row 2: Red in the Head / Leather Calot
$ws.Range("H2").Value = 413815.72
$ws.Range("I2").Value = 166618.33
$ws.Range("J2").Value = 599213.75
$ws.Range("K2").Value = 166618.33
$ws.Range("L2").Value = 599213.75
$ws.Range("M2").Value = -166506.33
$ws.Range("N2").Value = -599437.75

# Row 7: Tan Before the Ban / Leather
$ws.Range("H7").Value = 7835
$ws.Range("I7").Value = 5499
$ws.Range("J7").Value = 8224.333000000001
$ws.Range("K7").Value = 5499
$ws.Range("L7").Value = 8224.333000000001
$ws.Range("M7").Value = -5387
$ws.Range("N7").Value = -8448.333000000001

# Row 9: From the Sands to the Stage / Leather Himantes
$ws.Range("H9").Value = 1670.4
$ws.Range("I9").Value = 417.33334
$ws.Range("K9").Value = 417.33334
$ws.Range("M9").Value = -193.33334

# Row 40: Best Served Toad / Toad Leather
$ws.Range("H40").Value = 3831.6667
$ws.Range("I40").Value = 2774.25
$ws.Range("K40").Value = 2774.25
$ws.Range("M40").Value = -2638.25

# Row 126: Battered Books / Saiga Leather
$ws.Range("H126").Value = 7835
$ws.Range("I126").Value = 5499
$ws.Range("J126").Value = 8224.333000000001
$ws.Range("K126").Value = 16497
$ws.Range("L126").Value = 24672.999
$ws.Range("M126").Value = -14027
$ws.Range("N126").Value = -29612.999

$ws = $wb.Worksheets.Item("WVR")
# Row 122: Heavy Armoire / Dark Hempen Cloth
$ws.Range("H122").Value = 6105.1113
$ws.Range("I122").Value = 4675.2
$ws.Range("K122").Value = 14025.6
$ws.Range("M122").Value = -11575.6

# Row 136: Weaving the Envelope / Sarcenet Cloth
$ws.Range("H136").Value = 2344.2
$ws.Range("I136").Value = 2090.2307
$ws.Range("K136").Value = 6270.6921
$ws.Range("M136").Value = -3720.6921
